$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.716.10'
$ws.Range('E2').Value = '  +0.47%  '
$ws.Range('D3').Value = '1.849.58'
$ws.Range('D4').Value = "'1.030"
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = "'321.33"
$ws.Range('E5').Value = '  +0.45%  '
$ws.Range('D6').Value = "'1.028"
$ws.Range('E6').Value = '  +0.08%  '
$ws.Range('D7').Value = "'0.4386"
$ws.Range('E7').Value = '  +0.30%  '
$ws.Range('D8').Value = "'0.3789"
$ws.Range('E8').Value = '  +1.12%  '
$ws.Range('D9').Value = "'0.07395"
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('D10').Value = "'0.8823"
$ws.Range('E10').Value = '  +0.62%  '
$ws.Range('D11').Value = "'21.53"
$ws.Range('E11').Value = '  +0.23%  '
$ws.Range('D12').Value = '1.861.53'
$ws.Range('E12').Value = '  +0.29%  '
$ws.Range('D13').Value = "'5.502"
$ws.Range('E13').Value = '  +0.27%  '
$ws.Range('D14').Value = "'6.683"
$ws.Range('E14').Value = '  -0.03%  '
$ws.Range('D15').Value = "'0.07163"
$ws.Range('E15').Value = '  +0.20%  '
$ws.Range('D16').Value = "'84.87"
$ws.Range('E16').Value = '  +2.59%  '
$ws.Range('E17').Value = '  +0.22%  '
$ws.Range('D18').Value = "'0.000009063"
$ws.Range('E18').Value = '  +0.59%  '
$ws.Range('D19').Value = "'1.029"
$ws.Range('E19').Value = '  +0.22%  '
$ws.Range('D20').Value = "'15.46"
$ws.Range('E20').Value = '  +0.34%  '
$ws.Range('D21').Value = '27.726.98'
$ws.Range('E21').Value = '  +0.49%  '
$ws.Range('D22').Value = "'5.273"
$ws.Range('E22').Value = '  +0.25%  '
$ws.Range('D23').Value = "'11.24"
$ws.Range('E23').Value = '  +0.29%  '
$ws.Range('D24').Value = '2.096.43'
$ws.Range('E24').Value = '  +1.57%  '
$ws.Range('D25').Value = "'2.050"
$ws.Range('E25').Value = '  +6.26%  '
$ws.Range('D26').Value = "'158.42"
$ws.Range('E26').Value = '  +0.54%  '
$ws.Range('D27').Value = "'18.67"
$ws.Range('E27').Value = '  -0.28%  '
$ws.Range('D28').Value = "'1.991"
$ws.Range('E28').Value = '  +2.08%  '
$ws.Range('D29').Value = "'5.327"
$ws.Range('E29').Value = '  +1.41%  '
$ws.Range('D30').Value = "'117.67"
$ws.Range('E30').Value = '  +1.39%  '
$ws.Range('D31').Value = "'0.09061"
$ws.Range('E31').Value = '  -0.19%  '
$ws.Range('D32').Value = "'0.7746"
$ws.Range('E32').Value = '  +0.94%  '
$ws.Range('D33').Value = "'1.211"
$ws.Range('E33').Value = '  +0.44%  '
$ws.Range('D34').Value = "'2.994"
$ws.Range('E34').Value = '  +4.00%  '
$ws.Range('D35').Value = "'4.555"
$ws.Range('E35').Value = '  +1.21%  '
$ws.Range('D36').Value = "'1.030"
$ws.Range('E36').Value = '  +0.13%  '
$ws.Range('D37').Value = "'1.150"
$ws.Range('E37').Value = '  +0.50%  '
$ws.Range('D38').Value = "'0.01974"
$ws.Range('E38').Value = '  -0.22%  '
$ws.Range('D39').Value = "'0.05266"
$ws.Range('E39').Value = '  +0.09%  '
$ws.Range('D40').Value = "'2.845"
$ws.Range('E40').Value = '  +1.53%  '
$ws.Range('D41').Value = "'0.5178"
$ws.Range('E41').Value = '  +0.12%  '
$ws.Range('D42').Value = "'0.1669"
$ws.Range('E42').Value = '  -0.23%  '
$ws.Range('D43').Value = "'6.836"
$ws.Range('E43').Value = '  +2.34%  '
$ws.Range('D44').Value = "'8.723"
$ws.Range('E44').Value = '  +2.04%  '
$ws.Range('D45').Value = "'110.20"
$ws.Range('E45').Value = '  +1.25%  '
$ws.Range('E46').Value = '  +1.36%  '
$ws.Range('B47').Value = 'PaxDollar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D47').Value = "'1.031"
$ws.Range('E47').Value = '  +0.07%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = "'0.06573"
$ws.Range('E48').Value = '  +3.23%  '
$ws.Range('D49').Value = "'1.710"
$ws.Range('E49').Value = '  -0.30%  '
$ws.Range('D50').Value = "'0.4690"
$ws.Range('E50').Value = '  +0.80%  '
$ws.Range('D51').Value = "'1.883"
$ws.Range('E51').Value = '  -0.52%  '
